# Atualização automática de preços de eletricidade
# Updates the single data row (row 2) of the spot price table with the
# latest day's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (column A) - advance to the next day (stored as an Excel date serial)
$ws.Range("A2").Value = (Get-Date -Year 2026 -Month 2 -Day 11).Date

# Hourly prices (columns B..Z)
$ws.Range("B2").Value = 0.01
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -0.34
$ws.Range("E2").Value = -0.42
$ws.Range("F2").Value = -0.42
$ws.Range("G2").Value = -0.32
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.02
$ws.Range("J2").Value = 0.21
$ws.Range("K2").Value = 0.89
$ws.Range("L2").Value = 0.9
$ws.Range("M2").Value = 0.28
$ws.Range("N2").Value = 0.01
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0.02
$ws.Range("Q2").Value = 0.33
$ws.Range("R2").Value = 0.84
$ws.Range("S2").Value = 0.62
$ws.Range("T2").Value = 0.26
$ws.Range("U2").Value = 0.68
$ws.Range("V2").Value = 2.34
$ws.Range("W2").Value = 1.57
$ws.Range("X2").Value = 0.62
$ws.Range("Y2").Value = 0.11
$ws.Range("Z2").Value = 0.34

# Slot summary columns (AA..AG)
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 1.16
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 1.96
$ws.Range("AE2").Value = "16h-18h"
$ws.Range("AF2").Value = 0.73
$ws.Range("AG2").Value = "0h-23h"
